$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12)
$ws.Range("D2").Value = 10294
$ws.Range("E2").Value = 371
$ws.Range("F2").Value = 371
$ws.Range("G2").Value = 238
$ws.Range("H2").Value = 184
$ws.Range("I2").Value = 124
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 11484
$ws.Range("L2").Value = 6578
$ws.Range("M2").Value = 4906
$ws.Range("N2").Value = 4152
$ws.Range("O2").Value = 754
$ws.Range("P2").Value = 310
$ws.Range("Q2").Value = 423
$ws.Range("R2").Value = -526
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = 871
$ws.Range("U2").Value = -448
$ws.Range("V2").Value = 3896
$ws.Range("W2").Value = 3.6
$ws.Range("X2").Value = 1.79
$ws.Range("Y2").Value = 2.99
$ws.Range("Z2").Value = 1.59
$ws.Range("AA2").Value = 134.07
$ws.Range("AB2").Value = 1230.75
$ws.Range("AC2").Value = 200
$ws.Range("AD2").Value = 29.96
$ws.Range("AE2").Value = 6698
$ws.Range("AF2").Value = 0.9
$ws.Range("AG2").Value = 120
$ws.Range("AH2").Value = 2
$ws.Range("AI2").Value = 59.92
$ws.Range("AJ2").Value = 62000000

# Row 3 (2015/12)
$ws.Range("D3").Value = 10130
$ws.Range("E3").Value = 528
$ws.Range("F3").Value = 528
$ws.Range("G3").Value = 376
$ws.Range("H3").Value = 263
$ws.Range("I3").Value = 189
$ws.Range("J3").Value = 74
$ws.Range("K3").Value = 11575
$ws.Range("L3").Value = 6471
$ws.Range("M3").Value = 5104
$ws.Range("N3").Value = 4194
$ws.Range("O3").Value = 910
$ws.Range("P3").Value = 310
$ws.Range("Q3").Value = 968
$ws.Range("R3").Value = -413
$ws.Range("S3").Value = -488
$ws.Range("T3").Value = 637
$ws.Range("U3").Value = 331
$ws.Range("V3").Value = 3470
$ws.Range("W3").Value = 5.21
$ws.Range("X3").Value = 2.59
$ws.Range("Y3").Value = 4.52
$ws.Range("Z3").Value = 2.28
$ws.Range("AA3").Value = 126.79
$ws.Range("AB3").Value = 1248.84
$ws.Range("AC3").Value = 305
$ws.Range("AD3").Value = 15.89
$ws.Range("AE3").Value = 6764
$ws.Range("AF3").Value = 0.72
$ws.Range("AG3").Value = 125
$ws.Range("AH3").Value = 2.58
$ws.Range("AI3").Value = 41.04
$ws.Range("AJ3").Value = 62000000

# Row 4 (2016/12)
$ws.Range("D4").Value = 10761
$ws.Range("E4").Value = 597
$ws.Range("F4").Value = 597
$ws.Range("G4").Value = 572
$ws.Range("H4").Value = 461
$ws.Range("I4").Value = 384
$ws.Range("J4").Value = 78
$ws.Range("K4").Value = 12884
$ws.Range("L4").Value = 7342
$ws.Range("M4").Value = 5542
$ws.Range("N4").Value = 4491
$ws.Range("O4").Value = 1051
$ws.Range("P4").Value = 310
$ws.Range("Q4").Value = 906
$ws.Range("R4").Value = -875
$ws.Range("S4").Value = 200
$ws.Range("T4").Value = 884
$ws.Range("U4").Value = 22
$ws.Range("V4").Value = 3818
$ws.Range("W4").Value = 5.55
$ws.Range("X4").Value = 4.29
$ws.Range("Y4").Value = 8.84
$ws.Range("Z4").Value = 3.77
$ws.Range("AA4").Value = 132.49
$ws.Range("AB4").Value = 1350.26
$ws.Range("AC4").Value = 619
$ws.Range("AD4").Value = 7.24
$ws.Range("AE4").Value = 7243
$ws.Range("AF4").Value = 0.62
$ws.Range("AG4").Value = 125
$ws.Range("AH4").Value = 2.79
$ws.Range("AI4").Value = 20.2
$ws.Range("AJ4").Value = 62000000

# Row 5 (2017/12)
$ws.Range("D5").Value = 10282
$ws.Range("E5").Value = 382
$ws.Range("F5").Value = 382
$ws.Range("G5").Value = 297
$ws.Range("H5").Value = 262
$ws.Range("I5").Value = 218
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 11664
$ws.Range("L5").Value = 6076
$ws.Range("M5").Value = 5587
$ws.Range("N5").Value = 4572
$ws.Range("O5").Value = 1016
$ws.Range("P5").Value = 310
$ws.Range("Q5").Value = 933
$ws.Range("R5").Value = -331
$ws.Range("S5").Value = -817
$ws.Range("T5").Value = 479
$ws.Range("U5").Value = 455
$ws.Range("V5").Value = 3023
$ws.Range("W5").Value = 3.71
$ws.Range("X5").Value = 2.55
$ws.Range("Y5").Value = 4.8
$ws.Range("Z5").Value = 2.14
$ws.Range("AA5").Value = 108.75
$ws.Range("AB5").Value = 1404.45
$ws.Range("AC5").Value = 351
$ws.Range("AD5").Value = 10.66
$ws.Range("AE5").Value = 7374
$ws.Range("AF5").Value = 0.51
$ws.Range("AG5").Value = 115
$ws.Range("AH5").Value = 3.07
$ws.Range("AI5").Value = 32.77
$ws.Range("AJ5").Value = 62000000

# Row 6 (2018/12) -- note J6, O6 remain blank (not present before or after)
$ws.Range("D6").Value = 9954
$ws.Range("E6").Value = 422
$ws.Range("F6").Value = 422
$ws.Range("G6").Value = 314
$ws.Range("H6").Value = 239
$ws.Range("I6").Value = 174
$ws.Range("K6").Value = 10858
$ws.Range("L6").Value = 5150
$ws.Range("M6").Value = 5707
$ws.Range("N6").Value = 4656
$ws.Range("P6").Value = 310
$ws.Range("Q6").Value = 1199
$ws.Range("R6").Value = -328
$ws.Range("S6").Value = -947
$ws.Range("T6").Value = 265
$ws.Range("U6").Value = 934
$ws.Range("V6").Value = 2204
$ws.Range("W6").Value = 4.24
$ws.Range("X6").Value = 2.4
$ws.Range("Y6").Value = 3.77
$ws.Range("Z6").Value = 2.12
$ws.Range("AA6").Value = 90.25
$ws.Range("AB6").Value = 1449.84
$ws.Range("AC6").Value = 281
$ws.Range("AD6").Value = 13.8
$ws.Range("AE6").Value = 7509
$ws.Range("AF6").Value = 0.52
$ws.Range("AG6").Value = 115
$ws.Range("AH6").Value = 2.97
$ws.Range("AI6").Value = 40.95
$ws.Range("AJ6").Value = 62000000

# Rows 7, 8, 9 (2019/12(E), 2020/12(E), 2021/12(E)):
# all numeric data (columns D through AJ) is cleared, leaving only A, B, C labels.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
